# popAbwesenheitAnlegen.xlsx - "msz - Profil_Abwesenheiten_Normalfall_Anlage ok"
#
# The sheet documents an Appium/Android element locator profile. The
# Spinner/EditText locators were switched from @hint-based XPaths to
# @resource-id-based XPaths, the selected cell moved from F6 to E6, two
# columns were widened to fit the new (longer) locator strings, and the
# screenshot image was nudged/resized to match the new layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Widen columns C and E so the new, longer locator strings fit -------
# (column C: 26.21875 -> 45.109375 ; column E: 32.21875 -> 41.44140625,
#  expressed here in COM "characters" units)
$ws.Columns.Item(3).ColumnWidth = 44.333333333333336
$ws.Columns.Item(5).ColumnWidth = 40.666666666666664

# --- Update the Android element locators in row 2 -----------------------
# Column F ("//android.widget.Button[@text=\"Speichern\"]") is unchanged.
$ws.Range("C2").Value = '//android.widget.Spinner[@resource-id="from_date"]'
$ws.Range("D2").Value = '//android.widget.Spinner[@resource-id="till_date"]'
$ws.Range("E2").Value = '//android.widget.EditText[@resource-id="comment"]'

# --- Move the active selection from F6 to E6 -----------------------------
$ws.Range("E6").Select() | Out-Null

# --- Reposition/resize the screenshot picture to match the new layout ---
$shp = $ws.Shapes.Item(1)
$shp.Left = 0
$shp.Top = 106.8
$shp.Width = 1268.091496062992
$shp.Height = 569.1788188976378
